# Auto-generated edit script applying the diff's 127 cell changes
# across the 4 worksheets: 展览, 演出, 本地生活, 全部类型.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 42
$ws.Range("F3").Value = 8456
$ws.Range("F4").Value = 8456
$ws.Range("C6").Value = "【大会员提前抢】北京·ICOS内场-日本舞见鼻血姬"
$ws.Range("F6").Value = 6
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86902"
$ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202406/99SYO24h1717576009395.jpeg"
$ws.Range("F7").Value = 546
$ws.Range("F8").Value = 7442
$ws.Range("F9").Value = 7442
$ws.Range("F11").Value = 608
$ws.Range("F12").Value = 517
$ws.Range("F14").Value = 755
$ws.Range("F18").Value = 170
$ws.Range("F19").Value = 12333
$ws.Range("F20").Value = 108
$ws.Range("F21").Value = 14
$ws.Range("F22").Value = 2564
$ws.Range("F23").Value = 3675
$ws.Range("F24").Value = 57
$ws.Range("F26").Value = 2984
$ws.Range("F27").Value = 119
$ws.Range("F28").Value = 115
$ws.Range("F31").Value = 3381
$ws.Range("F34").Value = 1744
$ws.Range("F36").Value = 142
$ws.Range("F37").Value = 6095
$ws.Range("F39").Value = 1861
$ws.Range("F41").Value = 40
$ws.Range("F42").Value = 930
$ws.Range("F43").Value = 3
$ws.Range("F44").Value = 175
$ws.Range("F46").Value = 200
$ws.Range("F48").Value = 1115
$ws.Range("F49").Value = 1608
$ws.Range("F51").Value = 122
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 30
$ws.Range("F8").Value = 259
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 338
$ws.Range("F3").Value = 486
$ws.Range("F4").Value = 16
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 42
$ws.Range("F3").Value = 30
$ws.Range("F4").Value = 338
$ws.Range("F6").Value = 8456
$ws.Range("C7").Value = "丰台·【首家喜剧脱口秀】魔仙喜剧 l 舒哥专场 | 活的像一个笑话！"
$ws.Range("D7").Value = "丽泽天地购物中心 丽泽天地购物中心"
$ws.Range("E7").Value = "2024.07.20 17:30-07.20 20:30"
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 39.9
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=89211"
$ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202407/qaE2ZvuS1720778051895.jpeg"
$ws.Range("C8").Value = "北京·AINI二次元派对【免票展会】"
$ws.Range("D8").Value = "天竺镇裕翔路99号 北京欧陆时尚购物中心"
$ws.Range("E8").Value = "2024.07.20 16:00-07.21 19:00"
$ws.Range("F8").Value = 546
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=89134"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202407/lb2k8yDD1720678733848.jpeg"
$ws.Range("C9").Value = "北京·IDO动漫游戏嘉年华46th"
$ws.Range("D9").Value = "京沈路与天北路交汇处西北角 中国国际展览中心新馆"
$ws.Range("E9").Value = "2024.07.20 09:30-07.21 17:00"
$ws.Range("F9").Value = 7442
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=83716"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202405/9CAdQvG71716812495452.jpeg"
$ws.Range("C10").Value = "北京·原神同人嘉年华10th"
$ws.Range("D10").Value = "石景山路68号 北京首钢会展中心"
$ws.Range("E10").Value = "2024.07.20 09:00-07.21 17:00"
$ws.Range("F10").Value = 608
$ws.Range("G10").Value = 80
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=86012"
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/ulMhJXc61716260154833.jpeg"
$ws.Range("C11").Value = "北京·国乙同好嘉年华9th"
$ws.Range("F11").Value = 517
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=86011"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202405/AB4NkZsl1716272246698.jpeg"
$ws.Range("C12").Value = "北京·奶司的小人国娃展Nice Mini World"
$ws.Range("D12").Value = "永外高庄138号 大红门国际会展中心"
$ws.Range("E12").Value = "2024.07.20 10:30-07.20 17:00"
$ws.Range("F12").Value = 40
$ws.Range("G12").Value = 60
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=86952"
$ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202406/XD5Q2M561717658092443.jpeg"
$ws.Range("B13").Value = "'2024-07-21"
$ws.Range("C13").Value = "【大会员提前抢】北京·ICOS内场-青柳尊哉"
$ws.Range("D13").Value = "石景山路68号 北京首钢会展中心"
$ws.Range("E13").Value = "2024.07.21 09:00-07.21 17:00"
$ws.Range("F13").Value = 235
$ws.Range("G13").Value = 598
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86904"
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/WgYg2oTS1717576349209.jpeg"
$ws.Range("C14").Value = "北京·Summer Overture"
$ws.Range("D14").Value = "朝阳北路甲27号菁英梦谷·常营文创产业园南门B5座 WeShow Live 北京"
$ws.Range("E14").Value = "2024.07.21 12:00-07.21 19:00"
$ws.Range("F14").Value = 259
$ws.Range("G14").Value = 78
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=87481"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202406/dP7KKEIk1718608495643.png"
$ws.Range("C15").Value = "北京·世界名团首席系列—— 布达佩斯节日管弦乐团弦乐四重奏音乐会"
$ws.Range("D15").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$ws.Range("E15").Value = "2024.07.21 19:30-07.21 21:00"
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = 196
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=86891"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202405/wnG2Jyvg1717049167800.png"
$ws.Range("F17").Value = 170
$ws.Range("F19").Value = 12333
$ws.Range("F20").Value = 108
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 2564
$ws.Range("F24").Value = 3675
$ws.Range("F25").Value = 57
$ws.Range("F26").Value = 119
$ws.Range("F27").Value = 115
$ws.Range("F31").Value = 3381
$ws.Range("F33").Value = 1744
$ws.Range("F35").Value = 142
$ws.Range("F36").Value = 6095
$ws.Range("F39").Value = 1861
$ws.Range("F42").Value = 40
$ws.Range("F43").Value = 930
$ws.Range("F44").Value = 175
$ws.Range("F45").Value = 200
$ws.Range("F47").Value = 1115
$ws.Range("F49").Value = 1608
$ws.Range("F51").Value = 122
